# Update the "想去人数" (want-to-go count) values in column F
# for the "展览" and "全部类型" worksheets, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row => New Value, for sheet "展览"
$exhibitionUpdates = @{
    4  = 175
    5  = 4934
    9  = 541
    10 = 501
    13 = 1369
    14 = 3489
    16 = 129
    18 = 75
    19 = 2576
    24 = 42
    25 = 126
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row => New Value, for sheet "全部类型"
$allTypesUpdates = @{
    4  = 175
    6  = 4934
    10 = 541
    11 = 501
    14 = 1369
    15 = 3489
    17 = 129
    19 = 75
    20 = 2576
    25 = 42
    26 = 126
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
